$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: split the "[1 * 3]" label out of C1 into B1, and put a
#     multiplication sign "*" in C1 (matches row 13's existing layout
#     pattern of <size> <*> <size> across B/C/D). ---
$ws.Range("B1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = "*"

# --- Row 13: same fix applied to the second matrix-multiplication block. ---
$ws.Range("B13").Value = $ws.Range("C13").Value()
$ws.Range("C13").Value = "*"

# --- Update the visible selection / scroll position to match the new
#     view left after editing (was topLeftCell A4 / selection G13). ---
$ws.Range("G19").Select()
